$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09315466666666666
$ws.Range("N2").Value = 0.279464
$ws.Range("Q2").Value = 0.05943224261155555
$ws.Range("R2").Value = 0.534890183504
